# Fruta / hortaliza, semanal
# Update weekly price-report rows (2-14) on the active sheet: dates, variety,
# quality, volume and price columns are refreshed with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44371
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 7000
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7375
$ws.Range("P2").Value = 7375

# Row 3
$ws.Range("D3").Value = 44309
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("P3").Value = 8500

# Row 4
$ws.Range("D4").Value = 44414
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 6000
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 6500
$ws.Range("P4").Value = 6500

# Row 5
$ws.Range("D5").Value = 44253
$ws.Range("H5").Value = "Americana (o)"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 4000
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = 4250
$ws.Range("P5").Value = 4250

# Row 6 is unchanged

# Row 7
$ws.Range("D7").Value = 44636
$ws.Range("H7").Value = "Americana (o)"
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 8000
$ws.Range("L7").Value = 9000
$ws.Range("M7").Value = 8500
$ws.Range("P7").Value = 8500

# Row 8
$ws.Range("D8").Value = 44259
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 80

# Row 9
$ws.Range("D9").Value = 44539
$ws.Range("H9").Value = "Americana (o)"
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6750
$ws.Range("P9").Value = 6750

# Row 10
$ws.Range("D10").Value = 44497
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5500
$ws.Range("P10").Value = 5500

# Row 11
$ws.Range("D11").Value = 44410
$ws.Range("K11").Value = 5500
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5750
$ws.Range("P11").Value = 5750

# Row 12
$ws.Range("D12").Value = 44281
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5500
$ws.Range("P12").Value = 5500

# Row 13
$ws.Range("D13").Value = 44575
$ws.Range("H13").Value = "Sin especificar"

# Row 14
$ws.Range("D14").Value = 44263
$ws.Range("K14").Value = 7000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7500
$ws.Range("P14").Value = 7500
